# Auto-generated edit script applying the cryptos.xlsx diff.
# Refreshes the crypto price (col D) and 1h volume-change (col E) values,
# plus a few coin name/link swaps (rows 30/31, 41/42) and one full row
# replacement (row 51: Aave -> Algorand).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.970.70'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').Value = '2.549.57'
$ws.Range('E3').Value = '  +0.63%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = '''305.02'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('D6').Value = '''98.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.95%  '
$ws.Range('E7').Value = '  +0.98%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '''0.550'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').Value = '''36.94'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.66%  '
$ws.Range('E11').Value = '  +4.13%  '
$ws.Range('E12').Value = '  +1.60%  '
$ws.Range('D13').Value = '''7.65'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '2.940.33'
$ws.Range('E14').Value = '  +0.57%  '
$ws.Range('D15').Value = '2.548.85'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('E16').Value = '  +7.51%  '
$ws.Range('D17').Value = '''0.877'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').Value = '43.001.31'
$ws.Range('E18').Value = '  +0.45%  '
$ws.Range('D19').Value = '''13.86'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.30%  '
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').Value = '''71.96'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('D23').Value = '''255.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').Value = '''28.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.90%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '''10.26'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.77%  '
$ws.Range('D29').Value = '''37.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.41%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '''2.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = '''6.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.67%  '
$ws.Range('D32').Value = '''158.65'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.30%  '
$ws.Range('D33').Value = '''19.85'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +18.23%  '
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').Value = '''0.0803'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.48%  '
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('E37').Value = '  -4.55%  '
$ws.Range('E38').Value = '  +2.81%  '
$ws.Range('D39').Value = '''25.48'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.40%  '
$ws.Range('E40').Value = '  +0.57%  '
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').Value = '''3.45'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.14%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '''3.92'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.90%  '
$ws.Range('D43').Value = '''2.07'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +29.12%  '
$ws.Range('D44').Value = '2.107.40'
$ws.Range('E44').Value = '  +0.99%  '
$ws.Range('E45').Value = '  -0.99%  '
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').Value = '''86.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.33%  '
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('D49').Value = '2.798.95'
$ws.Range('E49').Value = '  +0.60%  '
$ws.Range('D50').Value = '''74.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +9.36%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').Value = '''0.193'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.52%  '
